# Update task row 2: title and status text changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Taraea actualizada"
$ws.Range("D2").Value = "En espera"
